# Atualizado por script em 03-01-2024 02:45
#
# This script:
#  1. Swaps the match-detail columns (F:V) between row 63 and row 65.
#  2. Swaps the match-detail columns (F:V) between row 66 and row 67.
#  3. Appends a new match record as row 184 (Valencia vs Villarreal).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A scratch row, well below the used data, to use as temporary holding
# space while swapping row contents.
$scratch = 1000

function Swap-Rows($rowA, $rowB) {
    $ws.Range("F$rowA`:V$rowA").Copy()
    $ws.Range("F$scratch`:V$scratch").PasteSpecial(-4104)

    $ws.Range("F$rowB`:V$rowB").Copy()
    $ws.Range("F$rowA`:V$rowA").PasteSpecial(-4104)

    $ws.Range("F$scratch`:V$scratch").Copy()
    $ws.Range("F$rowB`:V$rowB").PasteSpecial(-4104)

    $ws.Range("F$scratch`:V$scratch").ClearContents()
}

# 1) Rows 63 and 65 had their match data swapped.
Swap-Rows 63 65

# 2) Rows 66 and 67 had their match data swapped.
Swap-Rows 66 67

# 3) Append the new row 184 with match data (Valencia x Villarreal).
$ws.Range("A183:V183").Copy()
$ws.Range("A184:V184").PasteSpecial(-4122)

$ws.Cells.Item(184, 1).Value = 183
$ws.Cells.Item(184, 2).Value = "spain"
$ws.Cells.Item(184, 3).Value = "laliga"
$ws.Cells.Item(184, 4).Value = "2023-2024"
$ws.Cells.Item(184, 5).Value = 45293.89583333334
$ws.Cells.Item(184, 6).Value = "Valencia"
$ws.Cells.Item(184, 7).Value = 3
$ws.Cells.Item(184, 8).Value = "Villarreal"
$ws.Cells.Item(184, 9).Value = 1
$ws.Cells.Item(184, 10).Value = 2
$ws.Cells.Item(184, 11).Value = "17/12/2024 18:03"
$ws.Cells.Item(184, 12).Value = 2.23
$ws.Cells.Item(184, 13).Value = "02/01/2024 21:29"
$ws.Cells.Item(184, 14).Value = 3.56
$ws.Cells.Item(184, 15).Value = "17/12/2024 18:03"
$ws.Cells.Item(184, 16).Value = 3.43
$ws.Cells.Item(184, 17).Value = "02/01/2024 21:14"
$ws.Cells.Item(184, 18).Value = 3.49
$ws.Cells.Item(184, 19).Value = "17/12/2024 18:03"
$ws.Cells.Item(184, 20).Value = 3.49
$ws.Cells.Item(184, 21).Value = "02/01/2024 21:29"
$ws.Cells.Item(184, 22).Value = "https://www.betexplorer.com/football/spain/laliga/valencia-villarreal/bskYZGYp/"
